$wb = $excel.ActiveWorkbook

# Add new worksheet at the end and name it
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws.Name = "loginData"

# Header row
$ws.Range("A1").Value = "StartLoginTest"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "password"
$ws.Range("D1").Value = "runmode"

# runmode column first (matches original authoring order)
$ws.Range("D2").Value = "Y"
$ws.Range("D3").Value = "Y"
$ws.Range("D4").Value = "N"

# username / password columns
$ws.Range("B2").Value = "testname@gmail.com"
$ws.Range("C2").Value = "testing@1234"
$ws.Range("B3").Value = "testname@gmail.com"
$ws.Range("C3").Value = "testing@1234"
$ws.Range("B4").Value = "testname@gmail.com"
$ws.Range("C4").Value = "testing@1234"

# Borders for whole table
$ws.Range("A1:D4").Borders.LineStyle = 1

# Header fill
$ws.Range("A1:D1").Interior.Color = 65535

# Hyperlinks (font + style applied automatically on top of existing border)
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:testname@gmail.com")

# Apply hyperlink-style formatting to B4 first (then remove the link so only one
# hyperlink entry - covering B3:B4 - remains, while both cells keep the formatting)
$tmpLink = $ws.Hyperlinks.Add($ws.Range("B4"), "mailto:testname@gmail.com")
$tmpLink.Delete()
$ws.Hyperlinks.Add($ws.Range("B3:B4"), "mailto:testname@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "testname@gmail.com")
$ws.Range("B3:B4").Value = "testname@gmail.com"

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:testing@1234", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "Adobe@1234")
$ws.Range("C2").Value = "testing@1234"

$ws.Hyperlinks.Add($ws.Range("C3:C4"), "mailto:testing@1234", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "Adobe@1234")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:testing@1234", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "Adobe@1234")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:testing@1234", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "Adobe@1234")
$ws.Range("C3:C4").Value = "testing@1234"

# Column widths
$ws.Columns.Item(2).ColumnWidth = 24.5546875
$ws.Columns.Item(3).ColumnWidth = 12

# Selection / view
$ws.Range("C4").Select()
